$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Liz's standup row (row 2): fill in what she did, what she plans to do,
# and her obstacle (which also carries Word's automatic "last edit" bookmark).
$tbl.Cell(2, 2).Range.InsertAfter("Finished project page, Created editing in cart overview")
$tbl.Cell(2, 3).Range.InsertAfter("Finish editing in cart overview")
$tbl.Cell(2, 4).Range.InsertAfter("I’m trying to figure out how to get the admin logic to fit into the work already done. ~")

# Word leaves a "_GoBack" bookmark at the position of the most recent edit.
# Wrap it around the trailing marker character (wrapping a non-empty range
# is reliable), then strip the marker back out with Find & Replace so the
# bookmark collapses to sit right after the obstacle text.
$markerRange = $d.Content
$markerRange.Find.Execute("~") | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null

$cleanup = $d.Content
$cleanup.Find.Execute("~", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
